$d = $word.ActiveDocument

# 1. Delete the bullet paragraph: "El personal del hospital inicia sesión..."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*El personal del hospital inicia sesión en el sistema de gestión de donaciones de sangre.*") {
        $p.Range.Delete()
        break
    }
}

# 2. Merge runs for "1]: Si el donante ingresa..." text (replace with same text, single run)
$d.Content.Find.Execute("1]: Si el donante ingresa información incompleta o incorrecta, el sistema mostrará un mensaje de error y solicitará que se corrijan los campos necesarios.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "1]: Si el donante ingresa información incompleta o incorrecta, el sistema mostrará un mensaje de error y solicitará que se corrijan los campos necesarios.", 2)

# 3. Merge runs for "EX[2]: Si el donante intenta editar..." text
$d.Content.Find.Execute("EX[2]: Si el donante intenta editar o eliminar la información de un donante que no existe en el sistema, el sistema mostrará un mensaje de error",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "EX[2]: Si el donante intenta editar o eliminar la información de un donante que no existe en el sistema, el sistema mostrará un mensaje de error", 2)
